$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "CL" worksheet: add "urban residential" / "rural residential" /
#    "commercial" columns (B:D) with per-building-type lifetimes, and
#    relabel the header / row labels.
# ---------------------------------------------------------------------
$cl = $wb.Worksheets.Item("CL")

# New header strings, written in the order that reproduces the expected
# shared-string table ordering (B1, C1, D1 first, then the About note,
# then A1 last so "Building Component" is replaced by "Building
# Component (years)").
$cl.Range("B1").Value = "urban residential"
$cl.Range("C1").Value = "rural residential"
$cl.Range("D1").Value = "commercial"

# Match the bold/right-aligned header style used by B1 for the two new
# header cells.
$cl.Range("C1").Font.Bold = $true
$cl.Range("C1").HorizontalAlignment = -4152
$cl.Range("D1").Font.Bold = $true
$cl.Range("D1").HorizontalAlignment = -4152

# About!A29 note (added before CL!A1 so the shared-string index lines up)
$about = $wb.Worksheets.Item("About")
$about.Range("A29").Value = "For the U.S. model, we use the same component lifetimes across building types."

# Rename CL!A1 header last
$cl.Range("A1").Value = "Building Component (years)"

# Row labels A2:A7 - each now refers to the building component itself;
# values stay the same text, just confirm/re-set them explicitly.
$cl.Range("A2").Value = "heating"
$cl.Range("A3").Value = "cooling and ventilation"
$cl.Range("A4").Value = "envelope"
$cl.Range("A5").Value = "lighting"
$cl.Range("A6").Value = "appliances"
$cl.Range("A7").Value = "other component"

# New C/D formula columns mirroring column B (same lifetime for every
# building type in the U.S. model).
$cl.Range("C2").Formula = "=`$B2"
$cl.Range("D2").Formula = "=`$B2"
$cl.Range("C3:D7").Formula = "=`$B3"

# Match the number format/style used by column B (style index 13:
# integer number format) for the new formula cells.
$cl.Range("C2:D7").NumberFormat = "0"

# Column widths for the two new columns (closest achievable widths to
# 17.85546875 / 16.5703125 given this runtime's column-width rounding).
$cl.Columns.Item(3).ColumnWidth = 17
$cl.Columns.Item(4).ColumnWidth = 15.7

# Portrait page orientation for the CL sheet.
$cl.PageSetup.Orientation = 1

# Move the active cell to A2 on the CL sheet, then restore "About" as
# the active/selected sheet tab (matching the workbook's original active
# sheet).
$cl.Range("A2").Select()
$about.Select()
